$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with the collapsed checklist tuple-string
$ws.Range("A2").Value = "('Dark Ascension Checklist', ['Card', '(You can mark this card to represent a double-faced card in your library or hand.)', '" + [char]0x2610 + " Loyal Cathar {W}{W}', '" + [char]0x2610 + " Soul Seizer {3}{U}{U}', '" + [char]0x2610 + " Chosen of Markov {2}{B}', '" + [char]0x2610 + " Ravenous Demon {3}{B}{B}', '" + [char]0x2610 + " Afflicted Deserter {3}{R}', '" + [char]0x2610 + " Hinterland Hermit {1}{R}', '" + [char]0x2610 + " Mondronen Shaman {3}{R}', '" + [char]0x2610 + " Lambholt Elder {2}{G}', '" + [char]0x2610 + " Scorned Villager {1}{G}', '" + [char]0x2610 + " Wolfbitten Captive {G}', '" + [char]0x2610 + " Huntmaster of the Fells {2}{R}{G}', '" + [char]0x2610 + " Chalice of Life {3}', '" + [char]0x2610 + " Elbrus, the Binding Blade {7}'])"

# Update A3 with the Human token tuple-string
$ws.Range("A3").Value = "('Human', ['Token Creature " + [char]0x2014 + " Human', '1/1'])"

# Update A4 with the Sorin emblem tuple-string
$ws.Range("A4").Value = "('Sorin, Lord of Innistrad Emblem', ['Emblem " + [char]0x2014 + " Sorin', 'Creatures you control get +1/+0.'])"

# Update A5 with the Vampire token tuple-string
$ws.Range("A5").Value = "('Vampire', ['Token Creature " + [char]0x2014 + " Vampire', 'Lifelink', '1/1'])"

# Remove now-unused rows 6 through 27
$ws.Range("A6:A27").EntireRow.Delete()
